$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing row 1 (Yêu cầu / Nội dung) down to row 2, freeing up row 1 for the title/link
$ws.Rows("1:1").Insert()

# Row 1: title text that will become a hyperlink
$ws.Hyperlinks.Add($ws.Range("A1"), "https://www.freeprojectz.com/dfd/room-booking-system-dataflow-diagram", "", "", "https://www.freeprojectz.com/dfd/room-booking-system-dataflow-diagram")
$ws.Range("A1").Value = "Room Booking System Dataflow Diagram (DFD) FreeProjectz"

# Row 2: new header cell C2
$ws.Range("C2").Value = "Nhóm"

# Rows 3-8: management group names
$ws.Range("A3").Value = "Hotel Management"
$ws.Range("A4").Value = "Room Management"
$ws.Range("A5").Value = "Services Management"
$ws.Range("A6").Value = "Payment Management"
$ws.Range("A7").Value = "Booking Management"
$ws.Range("A8").Value = "Customer Management"

# Column widths to match bestFit columns in target (closest achievable values)
$ws.Columns("A").ColumnWidth = 54.416666666666664
$ws.Columns("B").ColumnWidth = 8.083333333333334

# Selection as shown in target sheetView
$ws.Range("C6").Select()
